$wb = $excel.ActiveWorkbook

# --- Step 1: remove the old "总计" sheet (frees its sheetId for reuse) ---
$zjOld = $wb.Worksheets.Item("总计")
$zjOld.Delete()

# --- Step 2: create "2022-Q1" by copying "2021-Q4" (same column layout), reusing the freed sheetId ---
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($null, $src)
$q1 = $wb.Worksheets.Item(6)
$q1.Name = "2022-Q1"

# --- Header row for 2022-Q1 ---
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# --- Data rows for 2022-Q1 ---
$q1.Range("B2").Value = "'501025"
$q1.Range("C2").Value = "'鹏华港股通中证香港银行投资指数（LOF）A"
$q1.Range("D2").Value = "'9.81"
$q1.Range("E2").Value = "'94.47"
$q1.Range("F2").Value = "'4.16"
$q1.Range("G2").Value = "'0.4081"
$q1.Range("H2").Value = 9
$q1.Range("B3").Value = "'010365"
$q1.Range("C3").Value = "'鹏华港股通中证香港银行投资指数（LOF）C"
$q1.Range("D3").Value = "'6.07"
$q1.Range("E3").Value = "'94.47"
$q1.Range("F3").Value = "'4.16"
$q1.Range("G3").Value = "'0.2525"
$q1.Range("H3").Value = 9
$q1.Range("B4").Value = "'001685"
$q1.Range("C4").Value = "'汇添富沪港深新价值股票"
$q1.Range("D4").Value = "'2.64"
$q1.Range("E4").Value = "'80.68"
$q1.Range("F4").Value = "'4.32"
$q1.Range("G4").Value = "'0.1140"
$q1.Range("H4").Value = 6
$q1.Range("B5").Value = "'005051"
$q1.Range("C5").Value = "'上投摩根标普港股通低波红利指数A"
$q1.Range("D5").Value = "'4.02"
$q1.Range("E5").Value = "'92.23"
$q1.Range("F5").Value = "'2.28"
$q1.Range("G5").Value = "'0.0917"
$q1.Range("H5").Value = 9
$q1.Range("B6").Value = "'006809"
$q1.Range("C6").Value = "'泰康港股通中证香港银行投资指数A"
$q1.Range("D6").Value = "'1.99"
$q1.Range("E6").Value = "'94.73"
$q1.Range("F6").Value = "'4.16"
$q1.Range("G6").Value = "'0.0828"
$q1.Range("H6").Value = 9
$q1.Range("B7").Value = "'005052"
$q1.Range("C7").Value = "'上投摩根标普港股通低波红利指数C"
$q1.Range("D7").Value = "'2.61"
$q1.Range("E7").Value = "'92.23"
$q1.Range("F7").Value = "'2.28"
$q1.Range("G7").Value = "'0.0595"
$q1.Range("H7").Value = 9
$q1.Range("B8").Value = "'501310"
$q1.Range("C8").Value = "'华宝标普沪港深中国增强价值指数（LOF）A"
$q1.Range("D8").Value = "'1.44"
$q1.Range("E8").Value = "'94.80"
$q1.Range("F8").Value = "'2.89"
$q1.Range("G8").Value = "'0.0416"
$q1.Range("H8").Value = 8
$q1.Range("B9").Value = "'006810"
$q1.Range("C9").Value = "'泰康港股通中证香港银行投资指数C"
$q1.Range("D9").Value = "'0.90"
$q1.Range("E9").Value = "'94.73"
$q1.Range("F9").Value = "'4.16"
$q1.Range("G9").Value = "'0.0374"
$q1.Range("H9").Value = 9
$q1.Range("B10").Value = "'007751"
$q1.Range("C10").Value = "'景顺长城中证沪港深红利成长低波动指数A"
$q1.Range("D10").Value = "'0.83"
$q1.Range("E10").Value = "'91.29"
$q1.Range("F10").Value = "'2.96"
$q1.Range("G10").Value = "'0.0246"
$q1.Range("H10").Value = 4
$q1.Range("B11").Value = "'004098"
$q1.Range("C11").Value = "'前海开源港股通股息率50强股票"
$q1.Range("D11").Value = "'0.34"
$q1.Range("E11").Value = "'88.92"
$q1.Range("F11").Value = "'3.52"
$q1.Range("G11").Value = "'0.0120"
$q1.Range("H11").Value = 5
$q1.Range("B12").Value = "'008407"
$q1.Range("C12").Value = "'恒生前海恒生沪深港通细分行业龙头指数A"
$q1.Range("D12").Value = "'0.37"
$q1.Range("E12").Value = "'93.40"
$q1.Range("F12").Value = "'2.68"
$q1.Range("G12").Value = "'0.0099"
$q1.Range("H12").Value = 9
$q1.Range("B13").Value = "'006106"
$q1.Range("C13").Value = "'景顺长城量化港股通股票"
$q1.Range("D13").Value = "'0.34"
$q1.Range("E13").Value = "'85.20"
$q1.Range("F13").Value = "'2.16"
$q1.Range("G13").Value = "'0.0073"
$q1.Range("H13").Value = 8
$q1.Range("B14").Value = "'001824"
$q1.Range("C14").Value = "'博时沪港深成长企业混合"
$q1.Range("D14").Value = "'0.15"
$q1.Range("E14").Value = "'93.41"
$q1.Range("F14").Value = "'2.87"
$q1.Range("G14").Value = "'0.0043"
$q1.Range("H14").Value = 10
$q1.Range("B15").Value = "'007397"
$q1.Range("C15").Value = "'华宝标普沪港深中国增强价值指数（LOF）C"
$q1.Range("D15").Value = "'0.09"
$q1.Range("E15").Value = "'94.80"
$q1.Range("F15").Value = "'2.89"
$q1.Range("G15").Value = "'0.0026"
$q1.Range("H15").Value = 8
$q1.Range("B16").Value = "'008408"
$q1.Range("C16").Value = "'恒生前海恒生沪深港通细分行业龙头指数C"
$q1.Range("D16").Value = "'0.08"
$q1.Range("E16").Value = "'93.40"
$q1.Range("F16").Value = "'2.68"
$q1.Range("G16").Value = "'0.0021"
$q1.Range("H16").Value = 9
$q1.Range("B17").Value = "'007760"
$q1.Range("C17").Value = "'景顺长城中证沪港深红利成长低波动指数C"
$q1.Range("D17").Value = "'0.06"
$q1.Range("E17").Value = "'91.29"
$q1.Range("F17").Value = "'2.96"
$q1.Range("G17").Value = "'0.0018"
$q1.Range("H17").Value = 4

# --- Step 3: create a brand-new "总计" sheet positioned right after "2022-Q1" (gets the next sheetId) ---
$zj = $wb.Worksheets.Add($null, $q1)
$zj.Name = "总计"
$zj.Outline.SummaryRow = 1
$zj.Outline.SummaryColumn = 1

# --- Copy header/row-label formatting (bold + border) from the 2022-Q1 sheet ---
$q1.Range("B1:D1").Copy()
$zj.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2:A7").Copy()
$zj.Range("A2:A7").PasteSpecial(-4122)

# --- Header row for 总计 ---
$zj.Range("B1").Value = "日期"
$zj.Range("C1").Value = "持有数量(只)"
$zj.Range("D1").Value = "持有市值(亿元)"

# --- Data rows for 总计 ---
$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 16
$zj.Range("D2").Value = 1.15
$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2021-Q4"
$zj.Range("C3").Value = 16
$zj.Range("D3").Value = 0.89
$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2021-Q3"
$zj.Range("C4").Value = 15
$zj.Range("D4").Value = 0.36
$zj.Range("A5").Value = 3
$zj.Range("B5").Value = "2021-Q2"
$zj.Range("C5").Value = 12
$zj.Range("D5").Value = 2
$zj.Range("A6").Value = 4
$zj.Range("B6").Value = "2021-Q1"
$zj.Range("C6").Value = 11
$zj.Range("D6").Value = 1.34
$zj.Range("A7").Value = 5
$zj.Range("B7").Value = "2020-Q4"
$zj.Range("C7").Value = 9
$zj.Range("D7").Value = 0.91
